$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(104).Insert()
$ws.Rows.Item(104).RowHeight = 25.5
$ws.Range("A104:B104").Merge()
$ws.Range("C104:G104").Merge()
$ws.Range("H104:K104").Merge()
$ws.Range("L104:M104").Merge()
$ws.Range("N104:O104").Merge()

Write-Host ("Row104 Height=" + $ws.Rows.Item(104).RowHeight)
Write-Host ("A104 mergearea=" + $ws.Range("A104").MergeArea.Address())
Write-Host ("C104 mergearea=" + $ws.Range("C104").MergeArea.Address())
Write-Host ("H104 mergearea=" + $ws.Range("H104").MergeArea.Address())
